# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) to the latest scraped figures.
#
# The Price column holds figures as *text* (e.g. "327.86", or
# "27.879.80" for the two-dot thousands-style values used by the big-cap
# coins). Excel's Range.Value setter auto-detects plain decimal-looking
# strings as numbers, so those are written with a leading apostrophe
# (the standard Excel "force text" input) to keep them as text, exactly
# like typing them into a cell by hand would. Values with two dots
# already fail numeric parsing and need no such prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.879.80"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "1.769.84"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'327.86"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.4486"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("D8").Value = "'0.3554"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "'0.07464"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "'42.10"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'1.098"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "'20.89"
$ws.Range("D14").Value = "'6.032"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "'7.224"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "1.767.16"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'93.20"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'0.00001059"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'0.06436"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'17.20"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").Value = "'5.775"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "27.912.28"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "'11.30"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'2.110"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "'161.92"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").Value = "'20.27"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "1.972.33"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("D30").Value = "'125.02"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "'1.097"
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("D32").Value = "'0.09172"
$ws.Range("D33").Value = "'5.584"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").Value = "'3.641"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'0.06115"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "'0.2098"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'4.966"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "'0.6297"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'1.181"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "'1.392"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").Value = "'13.28"
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "'0.5868"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "'122.57"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'1.953"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").Value = "'0.06906"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "'1.136"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "'73.03"
$ws.Range("E51").Value = "  +1.36%  "
